$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "hello",
    "ello",
    "hi",
    "hello",
    "what is order status",
    "what is order status",
    "hello",
    "what is order status",
    "what is order status",
    "what is status of last order",
    "when was it ?",
    "when was the last order delivered?"
)

$startRow = 99
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
